$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds a date value for every data row (2 through 321).
# The workbook was refreshed and every one of those dates moved forward by
# one day (45178 -> 45179), i.e. from 2023-09-05 to 2023-09-06.
for ($r = 2; $r -le 321; $r++) {
    $ws.Cells.Item($r, 3).Value = 45179
}
